$wb = $excel.ActiveWorkbook

# Rename the existing sheet and add a new one after it.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "forty_basket_currency"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Six_basket_currency"

# Populate the new sheet. Value assignment order matches the order the
# underlying shared strings were first introduced in the authored file.
$ws2.Range("A3").Value = "US dollar"
$ws2.Range("B3").Value = 9.1
$ws2.Range("C3").Value = 14.7
$ws2.Range("D3").Value = 11.6
$ws2.Range("E3").Value = 18.4

$ws2.Range("A1").Value = "Six_basket_currency"
$ws2.Range("B1").Value = "2015-16"
$ws2.Range("D1").Value = "2020-21 (P)"

$ws2.Range("A8").Value = "Hong Kong dollar"
$ws2.Range("B8").Value = 2.9
$ws2.Range("C8").Value = 4.8
$ws2.Range("D8").Value = 3.9
$ws2.Range("E8").Value = 4.6

$ws2.Range("A4").Value = "euro"
$ws2.Range("B4").Value = 11.4
$ws2.Range("C4").Value = 14
$ws2.Range("D4").Value = 11.6
$ws2.Range("E4").Value = 14.7

$ws2.Range("A5").Value = "Chinese yuan"
$ws2.Range("B5").Value = 10
$ws2.Range("C5").Value = 5
$ws2.Range("D5").Value = 12
$ws2.Range("E5").Value = 5.6

$ws2.Range("A6").Value = "British pound"
$ws2.Range("B6").Value = 2.2
$ws2.Range("C6").Value = 3.5
$ws2.Range("D6").Value = 2.2
$ws2.Range("E6").Value = 3.3

$ws2.Range("A7").Value = "Japanese yen"
$ws2.Range("B7").Value = 2.5
$ws2.Range("C7").Value = 2.3
$ws2.Range("D7").Value = 2.3
$ws2.Range("E7").Value = 1.7

$ws2.Range("B2").Value = "Trade-based Weight"
$ws2.Range("C2").Value = "Export-based Weight"
$ws2.Range("D2").Value = "Trade-based Weight"
$ws2.Range("E2").Value = "Export-based Weight"

# Selections / active cells on each sheet, matching the authored file.
[void]$ws1.Range("A44").Select()
[void]$ws2.Range("B22").Select()
[void]$ws2.Activate()

Write-Host "Done"
